$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: keyword changes from "passive income ideas" to "passive income"
$ws.Range("A8").Value = "passive income"

# Row 8 row height reverts to the default (12.8)
$ws.Rows.Item(8).RowHeight = 12.8

# Append three new keyword/appID rows
$ws.Range("A13").Value = "best bitcoin"
$ws.Range("B13").Value = "com.hamxa.shaynachim"

$ws.Range("A14").Value = "bitcoin course"
$ws.Range("B14").Value = "com.hamxa.shaynachim"

$ws.Range("A15").Value = "taxi"
$ws.Range("B15").Value = "com.singleton.strechy"

# Select the newly-added last cell, matching the author's final selection
$ws.Range("B15").Select()
